# Add 2022-Q4 data
#
# 1) Insert a new "2022-Q4" worksheet (with per-fund holding detail) right
#    before the existing "2022-Q3" worksheet.
# 2) Prepend a matching "2022-Q4" summary row at the top of the "总计" sheet,
#    shifting the previously-existing quarters' rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: "总计" (summary) sheet - shift rows down and insert the new quarter
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Give the newly-appended index cell (A9) the same style as the other index
# cells in column A (A2:A8 all carry the bold/centered/bordered style).
$ws1.Range("A8").Copy() | Out-Null
$ws1.Range("A9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$wb.Application.CutCopyMode = $false

# Shift the existing quarters' data down by one row (bottom-up, so the
# source values are never clobbered before they are copied). Column A (the
# running 0-based index) is left untouched - it already reads 0..6 on rows
# 2..8, and row 9 simply continues the sequence with 7.
$ws1.Range("B9").Value = "2020-Q4"
$ws1.Range("C9").Value = 11
$ws1.Range("D9").Value = 0.75
$ws1.Range("A9").Value = 7

$ws1.Range("B8").Value = "2021-Q1"
$ws1.Range("C8").Value = 12
$ws1.Range("D8").Value = 0.9

$ws1.Range("B7").Value = "2021-Q2"
$ws1.Range("C7").Value = 10
$ws1.Range("D7").Value = 1.06

$ws1.Range("B6").Value = "2021-Q3"
$ws1.Range("C6").Value = 11
$ws1.Range("D6").Value = 1.07

$ws1.Range("B5").Value = "2022-Q1"
$ws1.Range("C5").Value = 9
$ws1.Range("D5").Value = 0.78

$ws1.Range("B4").Value = "2022-Q2"
$ws1.Range("C4").Value = 8
$ws1.Range("D4").Value = 0.68

$ws1.Range("B3").Value = "2022-Q3"
$ws1.Range("C3").Value = 10
$ws1.Range("D3").Value = 0.5

# New quarter goes on top
$ws1.Range("B2").Value = "2022-Q4"
$ws1.Range("C2").Value = 7
$ws1.Range("D2").Value = 0.14

# ---------------------------------------------------------------------------
# Part 2: brand-new "2022-Q4" worksheet, inserted just before "2022-Q3"
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q4"

# Re-fetch the template sheet AFTER the insertion - sheet references taken
# before a Worksheets.Add() can end up stale once the collection shifts.
$template = $wb.Worksheets.Item("2022-Q3")

# --- Header row ---
$template.Range("B1:H1").Copy() | Out-Null
$newSheet.Range("B1:H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$wb.Application.CutCopyMode = $false

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# --- Column A (row index) style for the data rows ---
$template.Range("A2:A8").Copy() | Out-Null
$newSheet.Range("A2:A8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$wb.Application.CutCopyMode = $false

# --- Data rows ---
# Columns B,C,D,E,F,G hold text values in the source data (fund codes with
# leading zeros, and numbers formatted as fixed-point strings), so format
# them as Text first to stop Excel from silently converting them to numbers.
$newSheet.Range("B2:G8").NumberFormat = "@"

$data = @(
    @(0, "000179", "广发美国房地产指数（QDII）人民币A",          "1.60", "92.49", "2.79", "0.0446", 8),
    @(1, "000180", "广发美国房地产指数（QDII）美元A",            "1.60", "92.49", "2.79", "0.0446", 8),
    @(2, "160140", "南方道琼斯美国精选REIT指数（QDII-LOF）A",    "0.80", "92.31", "3.05", "0.0244", 7),
    @(3, "320017", "诺安全球收益不动产（QDII）",                  "0.24", "68.42", "5.15", "0.0124", 5),
    @(4, "160141", "南方道琼斯美国精选REIT指数（QDII-LOF）C",    "0.39", "92.31", "3.05", "0.0119", 7),
    @(5, "016278", "广发美国房地产指数（QDII）人民币C",          "0.01", "92.49", "2.79", "0.0003", 8),
    @(6, "016279", "广发美国房地产指数（QDII）美元C",            "0.01", "92.49", "2.79", "0.0003", 8)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# The temporary Text number-format was only needed to stop auto-conversion
# while typing the values in; now that the literal text is committed to the
# cells, put the formatting back to the sheet's normal (unstyled) look.
$newSheet.Range("B2:G8").Style = "Normal"

Write-Host "Added 2022-Q4 sheet and updated 总计 summary."
